$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "2025-04-29 09:44:14"
$ws.Range("B67").Value = 203
